# devoir2_results.xlsx -- "plots, sauvegarde et verbose, clean code"
#
# This script:
#   1. Refreshes the Gumbel (Q_gumbel / column D) return-period values on the
#      three existing "Return periods ..." sheets (the earlier values were
#      placeholders derived from the wrong parameterisation; the "inf" cells
#      for T_r(Q_2021) now resolve to finite numbers).
#   2. Appends two new worksheets summarising the goodness-of-fit tests that
#      were added to the analysis: "Kolmogorov-Smirnov" and "Chi-squared".

$wb = $excel.ActiveWorkbook

# ============================================================
# 1) Update the Gumbel return-period (column D) results on the
#    three existing 'Return periods ...' sheets.
# ============================================================
$ws = $wb.Worksheets.Item("Return periods Calen 20")
$ws.Range("D2").Value = 175.0626101605569
$ws.Range("D3").Value = 258.758871499615
$ws.Range("D4").Value = 340.9350656728547
$ws.Range("D5").Value = 422.9662722001634
$ws.Range("D6").Value = 679.3351850034544

$ws = $wb.Worksheets.Item("Return periods Calen 25")
$ws.Range("D2").Value = 188.8317091932907
$ws.Range("D3").Value = 286.5694763569031
$ws.Range("D4").Value = 382.5321586081328
$ws.Range("D5").Value = 478.3255290101483
$ws.Range("D6").Value = 264.485804983888

$ws = $wb.Worksheets.Item("Return periods Hydro 20")
$ws.Range("D2").Value = 176.7388488797245
$ws.Range("D3").Value = 259.4151098159642
$ws.Range("D4").Value = 340.5898285371615
$ws.Range("D5").Value = 421.6213265665889

# ============================================================
# 2) Add the two new goodness-of-fit result sheets at the end
#    of the workbook, styled to match the existing header row.
# ============================================================
$styleSrc = $wb.Worksheets.Item("Return periods Calen 20").Range("A1")

# ---- New sheet: Kolmogorov-Smirnov ----
$ks = $wb.Worksheets.Item($wb.Worksheets.Count)
$ks = $wb.Worksheets.Add($null, $ks)
$ks.Name = "Kolmogorov-Smirnov"

# Header row
$ks.Range("A1").Value = "Dataset"
$ks.Range("B1").Value = "Distribution"
$ks.Range("C1").Value = "D_max"
$ks.Range("D1").Value = "C_alpha"
$ks.Range("E1").Value = "alpha"
$ks.Range("F1").Value = "Result"

# Data rows
$ks.Range("A2").Value = "Calendar 2020"
$ks.Range("B2").Value = "Lognormal: Moments"
$ks.Range("C2").Value = 0.1084230343935367
$ks.Range("D2").Value = 0.1882500869537521
$ks.Range("E2").Value = 0.1
$ks.Range("F2").Value = "Accept H0"

$ks.Range("A3").Value = "Calendar 2020"
$ks.Range("B3").Value = "Lognormal: MaxLikelihood"
$ks.Range("C3").Value = 0.1155217666958316
$ks.Range("D3").Value = 0.1882500869537521
$ks.Range("E3").Value = 0.1
$ks.Range("F3").Value = "Accept H0"

$ks.Range("A4").Value = "Calendar 2020"
$ks.Range("B4").Value = "Gumbel"
$ks.Range("C4").Value = 0.1088799228167518
$ks.Range("D4").Value = 0.1882500869537521
$ks.Range("E4").Value = 0.1
$ks.Range("F4").Value = "Accept H0"

$ks.Range("A5").Value = "Hydrological 2020"
$ks.Range("B5").Value = "Lognormal: Moments"
$ks.Range("C5").Value = 0.09672734580030384
$ks.Range("D5").Value = 0.1882500869537521
$ks.Range("E5").Value = 0.1
$ks.Range("F5").Value = "Accept H0"

$ks.Range("A6").Value = "Hydrological 2020"
$ks.Range("B6").Value = "Lognormal: MaxLikelihood"
$ks.Range("C6").Value = 0.09347968874245371
$ks.Range("D6").Value = 0.1882500869537521
$ks.Range("E6").Value = 0.1
$ks.Range("F6").Value = "Accept H0"

$ks.Range("A7").Value = "Hydrological 2020"
$ks.Range("B7").Value = "Gumbel"
$ks.Range("C7").Value = 0.09881776540439746
$ks.Range("D7").Value = 0.1882500869537521
$ks.Range("E7").Value = 0.1
$ks.Range("F7").Value = "Accept H0"

# Apply header style (bold, border, center/top alignment) copied from an existing header cell
$styleSrc.Copy()
$ks.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- New sheet: Chi-squared ----
$chi = $wb.Worksheets.Item($wb.Worksheets.Count)
$chi = $wb.Worksheets.Add($null, $chi)
$chi.Name = "Chi-squared"

# Header row
$chi.Range("A1").Value = "Dataset"
$chi.Range("B1").Value = "Distribution"
$chi.Range("C1").Value = "k_classes"
$chi.Range("D1").Value = "Chi2_stat"
$chi.Range("E1").Value = "dof"
$chi.Range("F1").Value = "Critical_value"
$chi.Range("G1").Value = "alpha"
$chi.Range("H1").Value = "Result"

# Data rows
$chi.Range("A2").Value = "Calendar 2020"
$chi.Range("B2").Value = "Lognormal: Moments"
$chi.Range("C2").Value = 4
$chi.Range("D2").Value = 1.539002663778082
$chi.Range("E2").Value = 1
$chi.Range("F2").Value = 3.841458820694124
$chi.Range("G2").Value = 0.05
$chi.Range("H2").Value = "Accept H0"

$chi.Range("A3").Value = "Calendar 2020"
$chi.Range("B3").Value = "Lognormal: Moments"
$chi.Range("C3").Value = 5
$chi.Range("D3").Value = 4.46993198906311
$chi.Range("E3").Value = 2
$chi.Range("F3").Value = 5.991464547107979
$chi.Range("G3").Value = 0.05
$chi.Range("H3").Value = "Accept H0"

$chi.Range("A4").Value = "Calendar 2020"
$chi.Range("B4").Value = "Lognormal: MaxLikelihood"
$chi.Range("C4").Value = 4
$chi.Range("D4").Value = 1.778731921344894
$chi.Range("E4").Value = 1
$chi.Range("F4").Value = 3.841458820694124
$chi.Range("G4").Value = 0.05
$chi.Range("H4").Value = "Accept H0"

$chi.Range("A5").Value = "Calendar 2020"
$chi.Range("B5").Value = "Lognormal: MaxLikelihood"
$chi.Range("C5").Value = 5
$chi.Range("D5").Value = 5.101257355697612
$chi.Range("E5").Value = 2
$chi.Range("F5").Value = 5.991464547107979
$chi.Range("G5").Value = 0.05
$chi.Range("H5").Value = "Accept H0"

$chi.Range("A6").Value = "Calendar 2020"
$chi.Range("B6").Value = "Gumbel"
$chi.Range("C6").Value = 4
$chi.Range("D6").Value = 1.648332938984475
$chi.Range("E6").Value = 1
$chi.Range("F6").Value = 3.841458820694124
$chi.Range("G6").Value = 0.05
$chi.Range("H6").Value = "Accept H0"

$chi.Range("A7").Value = "Calendar 2020"
$chi.Range("B7").Value = "Gumbel"
$chi.Range("C7").Value = 5
$chi.Range("D7").Value = 4.564895208152135
$chi.Range("E7").Value = 2
$chi.Range("F7").Value = 5.991464547107979
$chi.Range("G7").Value = 0.05
$chi.Range("H7").Value = "Accept H0"

$chi.Range("A8").Value = "Hydrological 2020"
$chi.Range("B8").Value = "Lognormal: Moments"
$chi.Range("C8").Value = 4
$chi.Range("D8").Value = 1.674723409782879
$chi.Range("E8").Value = 1
$chi.Range("F8").Value = 3.841458820694124
$chi.Range("G8").Value = 0.05
$chi.Range("H8").Value = "Accept H0"

$chi.Range("A9").Value = "Hydrological 2020"
$chi.Range("B9").Value = "Lognormal: Moments"
$chi.Range("C9").Value = 5
$chi.Range("D9").Value = 7.254192607344461
$chi.Range("E9").Value = 2
$chi.Range("F9").Value = 5.991464547107979
$chi.Range("G9").Value = 0.05
$chi.Range("H9").Value = "Reject H0"

$chi.Range("A10").Value = "Hydrological 2020"
$chi.Range("B10").Value = "Lognormal: MaxLikelihood"
$chi.Range("C10").Value = 4
$chi.Range("D10").Value = 2.033592561050211
$chi.Range("E10").Value = 1
$chi.Range("F10").Value = 3.841458820694124
$chi.Range("G10").Value = 0.05
$chi.Range("H10").Value = "Accept H0"

$chi.Range("A11").Value = "Hydrological 2020"
$chi.Range("B11").Value = "Lognormal: MaxLikelihood"
$chi.Range("C11").Value = 5
$chi.Range("D11").Value = 7.956985019850496
$chi.Range("E11").Value = 2
$chi.Range("F11").Value = 5.991464547107979
$chi.Range("G11").Value = 0.05
$chi.Range("H11").Value = "Reject H0"

$chi.Range("A12").Value = "Hydrological 2020"
$chi.Range("B12").Value = "Gumbel"
$chi.Range("C12").Value = 4
$chi.Range("D12").Value = 1.76248375215539
$chi.Range("E12").Value = 1
$chi.Range("F12").Value = 3.841458820694124
$chi.Range("G12").Value = 0.05
$chi.Range("H12").Value = "Accept H0"

$chi.Range("A13").Value = "Hydrological 2020"
$chi.Range("B13").Value = "Gumbel"
$chi.Range("C13").Value = 5
$chi.Range("D13").Value = 7.388469580311943
$chi.Range("E13").Value = 2
$chi.Range("F13").Value = 5.991464547107979
$chi.Range("G13").Value = 0.05
$chi.Range("H13").Value = "Reject H0"

# Apply header style (bold, border, center/top alignment) copied from an existing header cell
$styleSrc.Copy()
$chi.Range("A1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the view the way it started (first sheet active, A1 selected)
$wb.Worksheets.Item(1).Activate()
$null = $wb.Worksheets.Item(1).Range("A1").Select()

